# CET 350 - Integrated Kevin's code
# Applies the narrative edits described in the commit diff using Word COM
# Find/Replace (wdReplaceAll) against the already-open ActiveDocument.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $result = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                                       $true, 1, $false, $replace, 2)
    if (-not $result) {
        Write-Host "WARNING: replace failed for:" $find
    }
}

# 1. Opening paragraph: add a topic sentence about the grandfather and change
#    "and fishing" -> "or fishing" in the list of pool activities.
Replace-Text `
    "Every summer, I would vacation with my grandparents to their condo in Siesta Key, Florida for a week or two.  I would spend the time by the pool, playing cards, reading a book, and fishing with my grandfather." `
    "The greatest man I’ve had the pleasure of knowing would be my grandfather.  Every summer, I would vacation with my grandparents to their condo in Siesta Key, Florida for a week or two.  I would spend the time by the pool, playing cards, reading a book, or fishing with my grandfather."

# 2. Biography paragraph: add "during the Great Depression", "school, and" ->
#    "school, but", "leave there" -> "leave home", and "first son" -> "first
#    of three sons".
Replace-Text `
    "Greene County, Pennsylvania.  He never finished high school, and would leave there to join the Army after the conclusion of World War II.  He would meet his wife, Mary Joan, and they would marry in 1949.  The couple would soon welcome their first son.  Shortly after," `
    "Greene County, Pennsylvania during the Great Depression.  He never finished high school, but would leave home to join the Army after the conclusion of World War II.  He would meet his wife, Mary Joan, and they would marry in 1949.  The couple would soon welcome their first of three sons.  Shortly after,"

# 3. "...starting as a dishwasher." -> "...starting as a dishwasher and working
#    his way up."
Replace-Text `
    "he would find work in a small bakery, starting as a dishwasher." `
    "he would find work in a small bakery, starting as a dishwasher and working his way up."

# 4. "...familiar with the establishment." -> "...familiar with the delicious
#    goodies that he made so popular."
Replace-Text `
    "those familiar with the establishment." `
    "those familiar with the delicious goodies that he made so popular."

# 5. "28-foot" -> "twenty-eight-foot"
Replace-Text `
    "fishing on his 28-foot fishing boat" `
    "fishing on his twenty-eight-foot fishing boat"

# 6. Rewrite of the generosity paragraph.
Replace-Text `
    "Aside from ensuring that my cousins and I would always be taken care of, he always had his hand out to his other family when they were in times of need.  I still recall a funny story of him buying my uncle’s farm so that uncle could buy an ugly, pink house for my aunt.  The financial aspects worked out in his favor eventually, as a bonus.  He left me with lesson that helping others when you have room on your shoulders is a virtue." `
    "He always his hands out to help his family in whatever way he could.  All his grandchildren, including myself, had access to a college fund.  Even outside of his immediate family, he made time for everyone.  I still recall a funny story of him buying my uncle’s farm so that uncle could buy an ugly, pink house for his wife. He left me with the lesson that helping others when you have room on your shoulders is a virtue."

# 7. "And much like his house," -> "Much like his house,"
Replace-Text `
    "we always joked he liked bleach more than water.  And much like his house," `
    "we always joked he liked bleach more than water.  Much like his house,"

$d.Saved = $false
